$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (sheet1): update Version + Date, insert a new "Jurisdiction" row ---
$meta = $wb.Worksheets.Item("Metadata")

# Insert a new row before row 11 (currently "Description") to make room for "Jurisdiction"
$meta.Rows.Item(11).Insert()

$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# Update Version and Date values (rows 3 and 8 are unaffected by the insert)
$meta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$meta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# --- Sheet "Elements" (sheet2): add the II-1 constraint text to the MaintainedEntity.typeId row ---
$elem = $wb.Worksheets.Item("Elements")
$elem.Range("AJ5").Value = "II-1:An II instance must have either a root or an nullFlavor. {root.exists() or nullFlavor.exists()}`n"
